# Add lab 2 data: fill in gender and age for subjects 7-12 (rows 8-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "m"
$ws.Range("C8").Value = 32

$ws.Range("B9").Value = "m"
$ws.Range("C9").Value = 25

$ws.Range("B10").Value = "f"
$ws.Range("C10").Value = 32

$ws.Range("B11").Value = "m"
$ws.Range("C11").Value = 28

$ws.Range("B12").Value = "f"
$ws.Range("C12").Value = 34

$ws.Range("B13").Value = "m"
$ws.Range("C13").Value = 63

# Set the active selection to match the end of the data-entry session
$ws.Range("B14").Select()
